$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = "Bar"
$ws.Range("C2").Value = "Bar"
$ws.Range("F2").Value = "     "
$ws.Range("H2").Value = "     "
$ws.Range("K2").Value = "     "
$ws.Range("L2").Value = "     "
$ws.Range("M2").Value = "     "
$ws.Range("N2").Value = "     "
$ws.Range("O2").Value = "     "
$ws.Range("C3").Value = "     "
$ws.Range("E3").Value = "     "
$ws.Range("G3").Value = "Bar"
$ws.Range("I3").Value = "     "
$ws.Range("K3").Value = "Bar"
$ws.Range("M3").Value = "Bar"
$ws.Range("B4").Value = "     "
$ws.Range("C4").Value = "     "
$ws.Range("I4").Value = "     "
$ws.Range("L4").Value = "Bar"
$ws.Range("M4").Value = "Bar"
$ws.Range("O4").Value = "Bar"
$ws.Range("E5").Value = "Bar"
$ws.Range("I5").Value = "Bar"
$ws.Range("M5").Value = "     "
$ws.Range("N5").Value = "Bar"
$ws.Range("O5").Value = "     "
$ws.Range("F6").Value = "Bar"
$ws.Range("H6").Value = "Bar"
$ws.Range("I6").Value = "Bar"
$ws.Range("J6").Value = "Bar"
$ws.Range("O6").Value = "Bar"
$ws.Range("B7").Value = "Server"
$ws.Range("C7").Value = "Bar"
$ws.Range("D7").Value = "Server"
$ws.Range("G7").Value = "Expo"
$ws.Range("I7").Value = "     "
$ws.Range("J7").Value = "Server"
$ws.Range("K7").Value = "Expo"
$ws.Range("M7").Value = "Server"
$ws.Range("O7").Value = "Server"
$ws.Range("B8").Value = "Server"
$ws.Range("D8").Value = "Expo"
$ws.Range("E8").Value = "Expo"
$ws.Range("F8").Value = "Expo"
$ws.Range("G8").Value = "     "
$ws.Range("H8").Value = "Expo"
$ws.Range("I8").Value = "Expo"
$ws.Range("K8").Value = "Expo"
$ws.Range("M8").Value = "Server"
$ws.Range("O8").Value = "Server"
$ws.Range("C9").Value = "     "
$ws.Range("J9").Value = "     "
$ws.Range("K9").Value = "     "
$ws.Range("N9").Value = "Server"
$ws.Range("E10").Value = "     "
$ws.Range("J10").Value = "Server"
$ws.Range("K10").Value = "     "
$ws.Range("L10").Value = "     "
$ws.Range("N10").Value = "     "
$ws.Range("O10").Value = "     "
$ws.Range("B11").Value = "     "
$ws.Range("I11").Value = "Server"
$ws.Range("N11").Value = "     "
$ws.Range("B12").Value = "     "
$ws.Range("D12").Value = "     "
$ws.Range("E12").Value = "     "
$ws.Range("G12").Value = "Server"
$ws.Range("L12").Value = "Server"
$ws.Range("M12").Value = "Server"
$ws.Range("B13").Value = "     "
$ws.Range("D13").Value = "     "
$ws.Range("F13").Value = "Server"
$ws.Range("G13").Value = "Server"
$ws.Range("H13").Value = "Server"
$ws.Range("J13").Value = "Server"
$ws.Range("K13").Value = "Server"
$ws.Range("L13").Value = "     "
$ws.Range("B14").Value = "     "
$ws.Range("C14").Value = "Server"
$ws.Range("F14").Value = "Server"
$ws.Range("J14").Value = "     "
$ws.Range("K14").Value = "Server"
$ws.Range("L14").Value = "Server"
$ws.Range("M14").Value = "     "
$ws.Range("D15").Value = "Server"
$ws.Range("E15").Value = "Server"
$ws.Range("K15").Value = "     "
$ws.Range("L15").Value = "Server"
$ws.Range("M15").Value = "     "
$ws.Range("N15").Value = "Server"
$ws.Range("B16").Value = "     "
$ws.Range("J16").Value = "     "
$ws.Range("K16").Value = "Server"
$ws.Range("L16").Value = "     "
$ws.Range("N16").Value = "Server"
$ws.Range("O16").Value = "Server"
$ws.Range("C17").Value = "Server"
$ws.Range("E17").Value = "     "
$ws.Range("G17").Value = "Server"
$ws.Range("H17").Value = "Server"
$ws.Range("M17").Value = "Server"
$ws.Range("N17").Value = "Server"
$ws.Range("B18").Value = "     "
$ws.Range("C18").Value = "     "
$ws.Range("E18").Value = "Server"
$ws.Range("G18").Value = "     "
$ws.Range("I18").Value = "Server"
$ws.Range("N18").Value = "Server"
$ws.Range("G19").Value = "Server"
$ws.Range("I19").Value = "Server"
$ws.Range("J19").Value = "Server"
$ws.Range("K19").Value = "     "
$ws.Range("L19").Value = "Server"
$ws.Range("C20").Value = "     "
$ws.Range("D20").Value = "Server"
$ws.Range("E20").Value = "Server"
$ws.Range("I20").Value = "     "
$ws.Range("N20").Value = "     "
$ws.Range("O20").Value = "     "
$ws.Range("B21").Value = "Expo"
$ws.Range("C21").Value = "Server"
$ws.Range("G21").Value = "Server"
$ws.Range("I21").Value = "Expo"
$ws.Range("J21").Value = "Expo"
$ws.Range("K21").Value = "Server"
$ws.Range("L21").Value = "Expo"
$ws.Range("N21").Value = "     "
$ws.Range("O21").Value = "Expo"
$ws.Range("G22").Value = "     "
$ws.Range("L22").Value = "     "
$ws.Range("B23").Value = "Server"
$ws.Range("G23").Value = "     "
$ws.Range("I23").Value = "     "
$ws.Range("K23").Value = "Server"
$ws.Range("O23").Value = "     "
$ws.Range("C24").Value = "Server"
$ws.Range("D24").Value = "     "
$ws.Range("E24").Value = "     "
$ws.Range("F24").Value = "     "
$ws.Range("G24").Value = "     "
$ws.Range("H24").Value = "     "
$ws.Range("M24").Value = "     "
$ws.Range("N24").Value = "     "
$ws.Range("B25").Value = "     "
$ws.Range("C25").Value = "     "
$ws.Range("E25").Value = "Server"
$ws.Range("F25").Value = "H/G"
$ws.Range("I25").Value = "     "
$ws.Range("K25").Value = "H/G"
$ws.Range("M25").Value = "     "
$ws.Range("O25").Value = "H/G"
$ws.Range("G26").Value = "     "
$ws.Range("J26").Value = "     "
$ws.Range("L26").Value = "     "
$ws.Range("N26").Value = "     "
$ws.Range("C27").Value = "H/G"
$ws.Range("E27").Value = "     "
$ws.Range("N27").Value = "H/G"
$ws.Range("O27").Value = "     "
$ws.Range("B28").Value = "H/G"
$ws.Range("C28").Value = "     "
$ws.Range("D28").Value = "     "
$ws.Range("E28").Value = "Runner"
$ws.Range("F28").Value = "     "
$ws.Range("I28").Value = "H/G"
$ws.Range("M28").Value = "     "
$ws.Range("O28").Value = "Runner"
$ws.Range("B29").Value = "Runner"
$ws.Range("D29").Value = "     "
$ws.Range("F29").Value = "     "
$ws.Range("I29").Value = "     "
$ws.Range("J29").Value = "     "
$ws.Range("K29").Value = "     "
$ws.Range("L29").Value = "     "
$ws.Range("M29").Value = "     "
$ws.Range("D30").Value = "H/G"
$ws.Range("L30").Value = "H/G"
$ws.Range("B31").Value = "     "
$ws.Range("G31").Value = "H/G"
$ws.Range("J31").Value = "H/G"
$ws.Range("M31").Value = "H/G"
$ws.Range("O31").Value = "     "
$ws.Range("C32").Value = "Expo"
$ws.Range("D32").Value = "Expo"
$ws.Range("E32").Value = "     "
$ws.Range("H32").Value = "     "
$ws.Range("I32").Value = "     "
$ws.Range("K32").Value = "     "
